$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Sheet3: add new row 6 ("question" label) ---
$ws3.Range("A6").Value = "question"
$ws3.Rows.Item(6).RowHeight = 21

# --- Sheet4: build the new data table ---
$ws4.Range("A1").Value = "question"
$ws4.Range("B1").Value = "question"
$ws4.Range("C1").Value = "answer1"
$ws4.Range("D1").Value = "answer2"
$ws4.Range("E1").Value = "answer3"

$ws4.Range("A2").Value = "dfaas"
$ws4.Range("C2").Value = "晚上好，不管白天还是黑夜都竭诚为你服务的人工智能机器人达尔文。我对答如流，无所不知，不知疲倦，不惧失败，帮你赚钱。你想了解我哪一点呢？"
$ws4.Range("D2").Value = "dasds"
$ws4.Range("E2").Value = "sadas"

$ws4.Range("B3").Value = "嗨喽 你知道我是谁吗？"
$ws4.Range("C3").Value = "gfartaeerawrdaf"
$ws4.Range("D3").Value = "我想你是我的朋友，也会成为我的忠实粉丝。"
$ws4.Range("E3").Value = "das"

$ws4.Range("B4").Value = "几天不见，变聪明了一点"
$ws4.Range("C4").Value = "safasd"
$ws4.Range("D4").Value = "sada"
$ws4.Range("E4").Value = "fsdad"

Write-Host "done"
